# Apply the 2023-09-26 04:41 UTC cryptos-list refresh (GitHub Actions update).
# Row 13/14 also swap coin identity (Polkadot <-> WrappedEther) in addition to
# updated Price / Volume(1h) figures throughout the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.354.06'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('E2').ClearFormats()

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.591.99'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E3').ClearFormats()

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('E4').ClearFormats()

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.41'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('E5').ClearFormats()

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.505'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E6').ClearFormats()

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('E7').ClearFormats()

# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E8').ClearFormats()

# Row 9
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.15%  '
$ws.Range('E9').ClearFormats()

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('E10').ClearFormats()

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('E11').ClearFormats()

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.815.98'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E12').ClearFormats()

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.627.07'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.76%  '
$ws.Range('E13').ClearFormats()

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('E14').ClearFormats()

# Row 15
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.03%  '
$ws.Range('E15').ClearFormats()

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.68'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('E16').ClearFormats()

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.356.14'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('E17').ClearFormats()

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('E18').ClearFormats()

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.76%  '
$ws.Range('E19').ClearFormats()

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '212.05'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.57%  '
$ws.Range('E20').ClearFormats()

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('E21').ClearFormats()

# Row 22
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('E22').ClearFormats()

# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('E23').ClearFormats()

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('E24').ClearFormats()

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.30%  '
$ws.Range('E25').ClearFormats()

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E26').ClearFormats()

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.07'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('E27').ClearFormats()

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.112'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('E28').ClearFormats()

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.21'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('E29').ClearFormats()

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0502'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E30').ClearFormats()

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('E31').ClearFormats()

# Row 32
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('E32').ClearFormats()

# Row 33
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('E33').ClearFormats()

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.341.69'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.31%  '
$ws.Range('E34').ClearFormats()

# Row 35
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.13%  '
$ws.Range('E35').ClearFormats()

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.601'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('E36').ClearFormats()

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E37').ClearFormats()

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0166'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('E38').ClearFormats()

# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -14.56%  '
$ws.Range('E39').ClearFormats()

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.818'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.25%  '
$ws.Range('E40').ClearFormats()

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.75'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +4.22%  '
$ws.Range('E41').ClearFormats()

# Row 43
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('E43').ClearFormats()

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.764'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('E44').ClearFormats()

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.728.45'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range('E45').ClearFormats()

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.64'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('E46').ClearFormats()

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '88.07'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.01%  '
$ws.Range('E47').ClearFormats()

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.50%  '
$ws.Range('E48').ClearFormats()

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.36%  '
$ws.Range('E49').ClearFormats()

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('E50').ClearFormats()

# Row 51
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.21%  '
$ws.Range('E51').ClearFormats()
